$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 92.4090904047999970
$ws.Range("C2").Value = [double]"1.85500908388E-9"

$ws.Range("B3").Value = 92409.0904047999940
$ws.Range("C3").Value = [double]"5.56502725164E-5"

$ws.Range("B4").Value = 341331.45722820982
$ws.Range("C4").Value = [double]"2.5573247980823871E-4"

$ws.Range("B5").Value = 6826.6291445641946
$ws.Range("C5").Value = [double]"5.1146495961647742E-6"
